$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price/volume data per the Aug 27 2023 GitHub Actions refresh.
# Column D (Price) values are forced to Text before assignment, then the cell style
# is reset to Normal, so purely numeric-looking strings (e.g. "1.010") are not
# auto-coerced into numbers by Excel's type inference, while keeping style index 0.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.426.28"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.18%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.673.09"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.17%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.010"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.80%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "221.04"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.52%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5366"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.51%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.010"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.74%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2670"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.36%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06416"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.54%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.08"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.39%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07859"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.86%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.583"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.34%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.682.04"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.29%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.901.08"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.11%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5646"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.03%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0₅8205"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.07%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "66.45"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.63%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "26.424.49"
$ws.Range("D18").Style = "Normal"

$ws.Range("E19").Value = "  +0.72%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.712"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.70%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "197.20"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.38%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.34"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.78%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.073"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.87%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.011"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.77%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "146.20"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.83%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1235"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.69%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.274"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.84%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "16.27"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.81%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.511"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.52%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05903"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.08%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.293"
$ws.Range("D31").Style = "Normal"

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.590"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.18%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.316"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.61%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.628"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.93%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9742"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.88%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.851"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.82%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.430"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.63%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5840"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.38%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01611"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.07%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.077.73"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +4.51%  "

$ws.Range("B41").Value = "FraxShare"
$ws.Range("C41").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.907"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.36%  "

$ws.Range("B42").Value = "TrustWalletToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.8668"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.52%  "

$ws.Range("E43").Value = "  +0.78%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "104.56"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.06%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.810.37"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.87%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "58.44"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.69%  "

$ws.Range("B47").Value = "BabyDogeCoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0₈107"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.47%  "

$ws.Range("B48").Value = "Frax"
$ws.Range("C48").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.014"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.26%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.4400"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.54%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.071"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.14%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05171"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.52%  "
